$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.362.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.32%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.653.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.60%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.57%  '

# Row 7
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$ws.Range("E8").Value = '  -0.88%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.654.36'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '

# Row 10
$ws.Range("E10").Value = '  -3.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.169'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.355'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '

# Row 13
$ws.Range("E13").Value = '  -1.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.142.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.16%  '

# Row 15
$ws.Range("E15").Value = '  -2.83%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '72.269.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.22%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.62%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.662.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.84%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '370.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.33%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.99%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

# Row 23
$ws.Range("E23").Value = '  -0.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.66%  '

# Row 25
$ws.Range("E25").Value = '  +0.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.85%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.797.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.97%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0970'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.17%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '494.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.94%  '

# Row 33
$ws.Range("E33").Value = '  -3.13%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.98%  '

# Row 35
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.83%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.07%  '

# Row 38
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '

# Row 40
$ws.Range("E40").Value = '  -2.97%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.54%  '

# Row 42
$ws.Range("E42").Value = '  -0.05%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.56%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.39%  '

# Row 45
$ws.Range("E45").Value = '  -0.90%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '155.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.96%  '

# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.51%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.556'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.96%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0756'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.42%  '
